$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 81; existing rows 81-94 shift down to 82-95.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with its new weekly record.
$ws.Range("A81").Value = 8
$ws.Range("B81").Value = "Terminal La Palmera de La Serena"
$ws.Range("C81").Value = "Coquimbo"
$ws.Range("D81").Value = 44511
$ws.Range("E81").Value = 4
$ws.Range("F81").Value = 100112040
$ws.Range("G81").Value = "Cilantro"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 3200
$ws.Range("K81").Value = 1300
$ws.Range("L81").Value = 1500
$ws.Range("M81").Value = 1400
$ws.Range("N81").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O81").Value = "Provincia del Elquí"
$ws.Range("P81").Value = 933
$ws.Range("Q81").Value = 1.5
$ws.Range("R81").Value = "Hortaliza"
